$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test description and test script name (shared strings)
$ws.Range("A5").Value = "Tax module should run fine on CICD"
$ws.Range("B5").Value = "Tax"

# Update the selected cell / view to B9
$ws.Range("B9").Select()
